$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price/volume columns), restoring a leading
# apostrophe so Excel keeps these numeric-looking strings as plain text,
# matching the original inlineStr cell contents.

$ws.Range("D2").Value = "'25.941.03"
$ws.Range("E2").Value = "'  -0.32%  "
$ws.Range("D3").Value = "'1.642.29"
$ws.Range("E3").Value = "'  +0.18%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("D5").Value = "'215.19"
$ws.Range("E5").Value = "'  -0.08%  "
$ws.Range("D6").Value = "'0.5046"
$ws.Range("E6").Value = "'  -0.13%  "
$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "'  -0.31%  "
$ws.Range("D8").Value = "'0.2572"
$ws.Range("E8").Value = "'  -0.25%  "
$ws.Range("D9").Value = "'0.06418"
$ws.Range("E9").Value = "'  +0.35%  "
$ws.Range("D10").Value = "'19.61"
$ws.Range("E10").Value = "'  +0.39%  "
$ws.Range("D11").Value = "'0.07776"
$ws.Range("D12").Value = "'4.267"
$ws.Range("E12").Value = "'  +0.54%  "
$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.641.26"
$ws.Range("E13").Value = "'  +0.00%  "
$ws.Range("B14").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'1.868.26"
$ws.Range("E14").Value = "'  +0.12%  "
$ws.Range("D15").Value = "'0.5436"
$ws.Range("E15").Value = "'  -0.27%  "
$ws.Range("D16").Value = "'0.0₅7933"
$ws.Range("E16").Value = "'  -0.08%  "
$ws.Range("D17").Value = "'64.55"
$ws.Range("E17").Value = "'  +1.44%  "
$ws.Range("D18").Value = "'25.971.13"
$ws.Range("E18").Value = "'  -0.27%  "
$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "'  -0.15%  "
$ws.Range("D20").Value = "'199.06"
$ws.Range("E20").Value = "'  -3.22%  "
$ws.Range("D21").Value = "'4.389"
$ws.Range("E21").Value = "'  +0.98%  "
$ws.Range("D22").Value = "'9.902"
$ws.Range("E22").Value = "'  -0.74%  "
$ws.Range("D23").Value = "'5.975"
$ws.Range("E23").Value = "'  -0.17%  "
$ws.Range("D24").Value = "'1.006"
$ws.Range("E24").Value = "'  -0.20%  "
$ws.Range("D25").Value = "'1.878"
$ws.Range("E25").Value = "'  -3.50%  "
$ws.Range("D26").Value = "'140.98"
$ws.Range("E26").Value = "'  -0.93%  "
$ws.Range("D27").Value = "'0.1136"
$ws.Range("E27").Value = "'  -1.89%  "
$ws.Range("B28").Value = "'EthereumClassic"
$ws.Range("C28").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'15.71"
$ws.Range("E28").Value = "'  -0.49%  "
$ws.Range("B29").Value = "'Cosmos"
$ws.Range("C29").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'6.811"
$ws.Range("E29").Value = "'  -0.78%  "
$ws.Range("D30").Value = "'1.242"
$ws.Range("E30").Value = "'  +0.38%  "
$ws.Range("D31").Value = "'0.04929"
$ws.Range("E31").Value = "'  -1.57%  "
$ws.Range("D32").Value = "'3.267"
$ws.Range("E32").Value = "'  -0.27%  "
$ws.Range("D33").Value = "'3.211"
$ws.Range("E33").Value = "'  +0.35%  "
$ws.Range("D34").Value = "'1.542"
$ws.Range("E34").Value = "'  +0.41%  "
$ws.Range("E35").Value = "'  +1.38%  "
$ws.Range("D36").Value = "'0.8933"
$ws.Range("E36").Value = "'  -1.59%  "
$ws.Range("D37").Value = "'2.607"
$ws.Range("E37").Value = "'  -1.66%  "
$ws.Range("D38").Value = "'1.148.77"
$ws.Range("E38").Value = "'  +2.17%  "
$ws.Range("D39").Value = "'0.5560"
$ws.Range("E39").Value = "'  -2.03%  "
$ws.Range("D40").Value = "'0.01573"
$ws.Range("E40").Value = "'  +0.72%  "
$ws.Range("D41").Value = "'1.005"
$ws.Range("E41").Value = "'  -0.17%  "
$ws.Range("D42").Value = "'5.717"
$ws.Range("E42").Value = "'  +1.40%  "
$ws.Range("D43").Value = "'0.8114"
$ws.Range("E43").Value = "'  -0.33%  "
$ws.Range("D44").Value = "'99.79"
$ws.Range("E44").Value = "'  +0.09%  "
$ws.Range("E45").Value = "'  +6.49%  "
$ws.Range("D46").Value = "'1.779.59"
$ws.Range("E46").Value = "'  +0.12%  "
$ws.Range("D47").Value = "'0.4531"
$ws.Range("E47").Value = "'  -0.04%  "
$ws.Range("E48").Value = "'  -0.40%  "
$ws.Range("D49").Value = "'54.75"
$ws.Range("E49").Value = "'  -0.35%  "
$ws.Range("D50").Value = "'0.05049"
$ws.Range("E50").Value = "'  +0.02%  "
$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "'  -0.05%  "
